$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. QDTAS1 (Sheet2) data updates
#    Row 2: email/password for user "onkark" updated
#    Row 3: email for user "Onkar" updated
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("QDTAS1")
$ws2.Cells.Item(2, 2).Value = "admin@gmail.com"
$ws2.Cells.Item(2, 3).Value = "Admin@123"
$ws2.Cells.Item(3, 2).Value = "tiselap800@irnini.com"

# ---------------------------------------------------------------------------
# 2. Sheet1 - append a new data row (row 6) for user "Shiv"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Item(6, 1).Value = "Shiv"
$ws1.Cells.Item(6, 2).Value = "misaxay339@fryshare.com"
$ws1.Cells.Item(6, 3).Value = 8745235686
$ws1.Cells.Item(6, 4).Value = "Shiv@12345"
$ws1.Cells.Item(6, 5).Value = "shivay@123"

# new hyperlinks for the password / newPass columns of the new row
$null = $ws1.Hyperlinks.Add($ws1.Range("D6"), "mailto:Shiv@12345")
$null = $ws1.Hyperlinks.Add($ws1.Range("E6"), "mailto:shivay@123")
$ws1.Range("D6").Style = "Hyperlink"
$ws1.Range("E6").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 3. Add a new, empty worksheet named "Sheet3" after QDTAS1
#    (add an extra sheet first and drop it so the internal sheetId counter
#     advances the same way it did in the authored workbook)
# ---------------------------------------------------------------------------
$tempSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheet3.Name = "Sheet3"
$null = $tempSheet.Delete()

# ---------------------------------------------------------------------------
# 4. Selections left the way the author left them
# ---------------------------------------------------------------------------
$null = $ws1.Activate()
$null = $ws1.Range("F8").Select()

$null = $ws2.Activate()
$null = $ws2.Range("B4").Select()
